$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.012014019070573
$ws.Range("D2").Value = 1.047310552161974
$ws.Range("E2").Value = 1.014121338360366
$ws.Range("F2").Value = 1.047650023309213
$ws.Range("I2").Value = 1.03816588614912
$ws.Range("J2").Value = 1.017259174851471
$ws.Range("K2").Value = 1.050073546448045
$ws.Range("L2").Value = 1.016980440022441
$ws.Range("M2").Value = 1.050412068711576
$ws.Range("N2").Value = 1.009901413314541

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.012962844597421
$ws.Range("D3").Value = 1.048062212629455
$ws.Range("E3").Value = 1.014924594023024
$ws.Range("F3").Value = 1.04863699319098
$ws.Range("I3").Value = 1.038347439040479
$ws.Range("J3").Value = 1.017841514589394
$ws.Range("K3").Value = 1.050637181457562
$ws.Range("L3").Value = 1.017588577923297
$ws.Range("M3").Value = 1.051210471632037
$ws.Range("N3").Value = 1.010093994106258

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.013576493200082
$ws.Range("D4").Value = 1.048544424930195
$ws.Range("E4").Value = 1.015444510948086
$ws.Range("F4").Value = 1.049271216130445
$ws.Range("I4").Value = 1.038461240550373
$ws.Range("J4").Value = 1.018217460067269
$ws.Range("K4").Value = 1.050997074577084
$ws.Range("L4").Value = 1.017981594509036
$ws.Range("M4").Value = 1.051722076158153
$ws.Range("N4").Value = 1.010218292242407

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.013834397612502
$ws.Range("D5").Value = 1.048746148980426
$ws.Range("E5").Value = 1.015663120583458
$ws.Range("F5").Value = 1.049536785914536
$ws.Range("I5").Value = 1.038508201746191
$ws.Range("J5").Value = 1.018375299233632
$ws.Range("K5").Value = 1.051147218976202
$ws.Range("L5").Value = 1.018146700804539
$ws.Range("M5").Value = 1.051935952449633
$ws.Range("N5").Value = 1.010270471642466

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.013877696576633
$ws.Range("D6").Value = 1.048779960765535
$ws.Range("E6").Value = 1.015699828218004
$ws.Range("F6").Value = 1.049581314198552
$ws.Range("I6").Value = 1.038516035045194
$ws.Range("J6").Value = 1.018401788888003
$ws.Range("K6").Value = 1.051172361121985
$ws.Range("L6").Value = 1.018174415943626
$ws.Range("M6").Value = 1.051971792672432
$ws.Range("N6").Value = 1.010279228356609

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.013579939620762
$ws.Range("D7").Value = 1.048547124300581
$ws.Range("E7").Value = 1.015447431878382
$ws.Range("F7").Value = 1.04927476884419
$ws.Range("I7").Value = 1.03846187151128
$ws.Range("J7").Value = 1.01821956994141
$ws.Range("K7").Value = 1.050999085353324
$ws.Range("L7").Value = 1.017983801131617
$ws.Range("M7").Value = 1.051724938709987
$ws.Range("N7").Value = 1.010218989763242

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.012334741982418
$ws.Range("D8").Value = 1.047565440275052
$ws.Range("E8").Value = 1.014392769328016
$ws.Range("F8").Value = 1.047984487399638
$ws.Range("I8").Value = 1.038228003180461
$ws.Range("J8").Value = 1.017456158465643
$ws.Range("K8").Value = 1.0502650258809
$ws.Range("L8").Value = 1.017186063952677
$ws.Range("M8").Value = 1.050682930350899
$ws.Range("N8").Value = 1.009966561875403

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.010138226728848
$ws.Range("D9").Value = 1.045803770739517
$ws.Range("E9").Value = 1.012535561363039
$ws.Range("F9").Value = 1.04567709467411
$ws.Range("I9").Value = 1.037787792154339
$ws.Range("J9").Value = 1.016104315807701
$ws.Range("K9").Value = 1.048934695558413
$ws.Range("L9").Value = 1.015776632370311
$ws.Range("M9").Value = 1.048808424343099
$ws.Range("N9").Value = 1.009519352860898

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.00867234336161
$ws.Range("D10").Value = 1.044608003329458
$ws.Range("E10").Value = 1.011298309754793
$ws.Range("F10").Value = 1.044116183192049
$ws.Range("I10").Value = 1.037475474573586
$ws.Range("J10").Value = 1.015198674223371
$ws.Range("K10").Value = 1.048023125838685
$ws.Range("L10").Value = 1.014834550077651
$ws.Range("M10").Value = 1.04753303204222
$ws.Range("N10").Value = 1.00921961278282

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.008037236735509
$ws.Range("D11").Value = 1.044085177179233
$ws.Range("E11").Value = 1.01076278632216
$ws.Range("F11").Value = 1.043434923770157
$ws.Range("I11").Value = 1.037335777472135
$ws.Range("J11").Value = 1.014805480542724
$ws.Range("K11").Value = 1.047622564169062
$ws.Range("L11").Value = 1.014426041098698
$ws.Range("M11").Value = 1.046974681905786
$ws.Range("N11").Value = 1.009089444489583

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.007801274373927
$ws.Range("D12").Value = 1.0438902180045
$ws.Range("E12").Value = 1.010563902013483
$ws.Range("F12").Value = 1.043181066195385
$ws.Range("I12").Value = 1.037283218020801
$ws.Range("J12").Value = 1.014659274268896
$ws.Range("K12").Value = 1.047472900127002
$ws.Range("L12").Value = 1.014274215811134
$ws.Range("M12").Value = 1.046766370202106
$ws.Range("N12").Value = 1.009041037402754

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.007851891612148
$ws.Range("D13").Value = 1.0439320717141
$ws.Range("E13").Value = 1.01060656187162
$ws.Range("F13").Value = 1.043235556049759
$ws.Range("I13").Value = 1.037294522499703
$ws.Range("J13").Value = 1.014690643093692
$ws.Range("K13").Value = 1.047505043275621
$ws.Range("L13").Value = 1.014306786774663
$ws.Range("M13").Value = 1.046811095197928
$ws.Range("N13").Value = 1.009051423454917

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.008017733144903
$ws.Range("D14").Value = 1.044069077258422
$ws.Range("E14").Value = 1.010746345807382
$ws.Range("F14").Value = 1.043413956297196
$ws.Range("I14").Value = 1.037331446550702
$ws.Range("J14").Value = 1.014793398279905
$ws.Range("K14").Value = 1.047610210791949
$ws.Range("L14").Value = 1.014413492942607
$ws.Range("M14").Value = 1.0469574814717
$ws.Range("N14").Value = 1.009085444302395

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.008119906234854
$ws.Range("D15").Value = 1.044153390434237
$ws.Range("E15").Value = 1.01083247568826
$ws.Range("F15").Value = 1.04352376757331
$ws.Range("I15").Value = 1.037354107947532
$ws.Range("J15").Value = 1.014856688395843
$ws.Range("K15").Value = 1.047674891697329
$ws.Range("L15").Value = 1.014479226636865
$ws.Range("M15").Value = 1.04704755356908
$ws.Range("N15").Value = 1.009106398148306

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.008714485515535
$ws.Range("D16").Value = 1.044642595213887
$ws.Range("E16").Value = 1.011333855263576
$ws.Range("F16").Value = 1.044161282902347
$ws.Range("I16").Value = 1.037484651860698
$ws.Range("J16").Value = 1.015224747208895
$ws.Range("K16").Value = 1.048049586691059
$ws.Range("L16").Value = 1.014861649268089
$ws.Range("M16").Value = 1.047569959483895
$ws.Range("N16").Value = 1.009228243653149

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.009087350504147
$ws.Range("D17").Value = 1.044948108612615
$ws.Range("E17").Value = 1.011648415394646
$ws.Range("F17").Value = 1.044559740548453
$ws.Range("I17").Value = 1.037565344551753
$ws.Range("J17").Value = 1.015455341288029
$ws.Range("K17").Value = 1.048283058353627
$ws.Range("L17").Value = 1.015101377619986
$ws.Range("M17").Value = 1.047896018925785
$ws.Range("N17").Value = 1.009304572784559

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.009304800696941
$ws.Range("D18").Value = 1.045125822030959
$ws.Range("E18").Value = 1.011831913666372
$ws.Range("F18").Value = 1.044791635606175
$ws.Range("I18").Value = 1.037611980726121
$ws.Range("J18").Value = 1.015589742074005
$ws.Range("K18").Value = 1.048418674333875
$ws.Range("L18").Value = 1.015241151013171
$ws.Range("M18").Value = 1.048085615605919
$ws.Range("N18").Value = 1.009349057712628

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.009378939545727
$ws.Range("D19").Value = 1.045186335026657
$ws.Range("E19").Value = 1.011894485319527
$ws.Range("F19").Value = 1.044870617872996
$ws.Range("I19").Value = 1.037627809416729
$ws.Range("J19").Value = 1.015635552153558
$ws.Range("K19").Value = 1.048464820176341
$ws.Range("L19").Value = 1.015288800598901
$ws.Range("M19").Value = 1.04815016337881
$ws.Range("N19").Value = 1.009364219717598

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.009047349306285
$ws.Range("D20").Value = 1.044915380323142
$ws.Range("E20").Value = 1.011614663926481
$ws.Range("F20").Value = 1.044517043449035
$ws.Range("I20").Value = 1.037556731514577
$ws.Range("J20").Value = 1.015430611153156
$ws.Range("K20").Value = 1.048258067394657
$ws.Range("L20").Value = 1.015075662835726
$ws.Range("M20").Value = 1.047861096677988
$ws.Range("N20").Value = 1.009296387171971

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.007968898478562
$ws.Range("D21").Value = 1.044028753468045
$ws.Range("E21").Value = 1.010705182007716
$ws.Range("F21").Value = 1.043361444162602
$ws.Range("I21").Value = 1.03732059182153
$ws.Range("J21").Value = 1.014763143759892
$ws.Range("K21").Value = 1.047579265800814
$ws.Range("L21").Value = 1.014382073030183
$ws.Range("M21").Value = 1.046914399635389
$ws.Range("N21").Value = 1.009075427577348

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.007290512622459
$ws.Range("D22").Value = 1.043466909099633
$ws.Range("E22").Value = 1.010133545515785
$ws.Range("F22").Value = 1.042630201001971
$ws.Range("I22").Value = 1.037168246818763
$ws.Range("J22").Value = 1.014342574431033
$ws.Range("K22").Value = 1.047147399013056
$ws.Range("L22").Value = 1.013945483220793
$ws.Range("M22").Value = 1.046313876325283
$ws.Range("N22").Value = 1.00893617301605

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.007650168229927
$ws.Range("D23").Value = 1.043765169090006
$ws.Range("E23").Value = 1.01043656252152
$ws.Range("F23").Value = 1.043018289742023
$ws.Range("I23").Value = 1.037249374883186
$ws.Range("J23").Value = 1.014565611968525
$ws.Range("K23").Value = 1.047376820837833
$ws.Range("L23").Value = 1.014176975218021
$ws.Range("M23").Value = 1.046632727127242
$ws.Range("N23").Value = 1.009010025589895

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.009065424235686
$ws.Range("D24").Value = 1.044930170333811
$ws.Range("E24").Value = 1.011629914698195
$ws.Range("F24").Value = 1.044536338031884
$ws.Range("I24").Value = 1.037560624705306
$ws.Range("J24").Value = 1.015441785948053
$ws.Range("K24").Value = 1.048269361475779
$ws.Range("L24").Value = 1.015087282413089
$ws.Range("M24").Value = 1.047876878355596
$ws.Range("N24").Value = 1.009300086010887

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.010706351208077
$ws.Range("D25").Value = 1.046262967234636
$ws.Range("E25").Value = 1.013015541463022
$ws.Range("F25").Value = 1.046277605769122
$ws.Range("I25").Value = 1.037904922164714
$ws.Range("J25").Value = 1.016454579721014
$ws.Range("K25").Value = 1.049282973404613
$ws.Range("L25").Value = 1.016141440665199
$ws.Range("M25").Value = 1.049297566934042
$ws.Range("N25").Value = 1.009635250092434
